$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '65.966.71'
$ws.Range('E2').Value = '  -4.95%  '

$ws.Range('D3').Value = '3.309.56'
$ws.Range('E3').Value = '  -6.28%  '

$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.28%  '

$ws.Range('D5').Value = '558.91'
$ws.Range('E5').Value = '  -4.59%  '

$ws.Range('D6').Value = '180.92'
$ws.Range('E6').Value = '  -6.68%  '

$ws.Range('D7').Value = '0.999'
$ws.Range('E7').Value = '  -0.09%  '

$ws.Range('D8').Value = '0.591'
$ws.Range('E8').Value = '  -2.70%  '

$ws.Range('D9').Value = '3.306.83'
$ws.Range('E9').Value = '  -5.99%  '

$ws.Range('D10').Value = '0.189'
$ws.Range('E10').Value = '  -7.36%  '

$ws.Range('D11').Value = '0.589'
$ws.Range('E11').Value = '  -5.12%  '

$ws.Range('D12').Value = '47.78'
$ws.Range('E12').Value = '  -8.74%  '

$ws.Range('D13').Value = '0.0000266'
$ws.Range('E13').Value = '  -7.33%  '

$ws.Range('B14').Value = 'BitcoinCash'
$ws.Range('C14').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D14').Value = '636.44'
$ws.Range('E14').Value = '  -1.61%  '

$ws.Range('B15').Value = 'Polkadot'
$ws.Range('C15').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D15').Value = '8.58'
$ws.Range('E15').Value = '  -6.89%  '

$ws.Range('D16').Value = '3.829.04'
$ws.Range('E16').Value = '  -6.65%  '

$ws.Range('D17').Value = '18.11'
$ws.Range('E17').Value = '  -1.26%  '

$ws.Range('D18').Value = '65.910.44'
$ws.Range('E18').Value = '  -5.15%  '

$ws.Range('E19').Value = '  -3.84%  '

$ws.Range('D20').Value = '3.297.48'
$ws.Range('E20').Value = '  -6.98%  '

$ws.Range('D21').Value = '11.46'
$ws.Range('E21').Value = '  -8.50%  '

$ws.Range('D22').Value = '0.908'
$ws.Range('E22').Value = '  -5.12%  '

$ws.Range('D23').Value = '17.69'
$ws.Range('E23').Value = '  -2.31%  '

$ws.Range('D24').Value = '107.90'
$ws.Range('E24').Value = '  +5.87%  '

$ws.Range('D25').Value = '5.06'
$ws.Range('E25').Value = '  -7.69%  '

$ws.Range('D26').Value = '4.02'
$ws.Range('E26').Value = '  -8.02%  '

$ws.Range('E27').Value = '  -0.46%  '

$ws.Range('D28').Value = '2.69'
$ws.Range('E28').Value = '  -7.56%  '

$ws.Range('D29').Value = '9.55'
$ws.Range('E29').Value = '  -5.71%  '

$ws.Range('D30').Value = '8.80'
$ws.Range('E30').Value = '  -7.41%  '

$ws.Range('D31').Value = '30.75'
$ws.Range('E31').Value = '  -6.80%  '

$ws.Range('D32').Value = '3.99'
$ws.Range('E32').Value = '  -2.78%  '

$ws.Range('D33').Value = '6.34'
$ws.Range('E33').Value = '  -6.14%  '

$ws.Range('D34').Value = '11.09'
$ws.Range('E34').Value = '  -5.16%  '

$ws.Range('D35').Value = '551.62'
$ws.Range('E35').Value = '  +8.22%  '

$ws.Range('D36').Value = '0.106'
$ws.Range('E36').Value = '  -3.75%  '

$ws.Range('D37').Value = '3.738.62'
$ws.Range('E37').Value = '  +0.35%  '

$ws.Range('D38').Value = '1.00'
$ws.Range('E38').Value = '  +0.11%  '

$ws.Range('D39').Value = '57.21'
$ws.Range('E39').Value = '  -6.98%  '

$ws.Range('D40').Value = '3.48'
$ws.Range('E40').Value = '  -2.93%  '

$ws.Range('D41').Value = '2.75'
$ws.Range('E41').Value = '  -6.84%  '

$ws.Range('D42').Value = '0.0₃0716'
$ws.Range('E42').Value = '  -10.83%  '

$ws.Range('E43').Value = '  +24.89%  '

$ws.Range('D44').Value = '0.127'
$ws.Range('E44').Value = '  -5.12%  '

$ws.Range('D45').Value = '0.343'
$ws.Range('E45').Value = '  -7.18%  '

$ws.Range('D46').Value = '32.18'
$ws.Range('E46').Value = '  -7.08%  '

$ws.Range('D47').Value = '0.0414'
$ws.Range('E47').Value = '  -7.13%  '

$ws.Range('D48').Value = '3.23'
$ws.Range('E48').Value = '  -5.24%  '

$ws.Range('E49').Value = '  -7.88%  '

$ws.Range('D50').Value = '0.130'
$ws.Range('E50').Value = '  -4.65%  '

$ws.Range('E51').Value = '  -0.24%  '
